$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.930.96"
$ws.Range("E2").Value = "  -3.66%  "

$ws.Range("D3").Value = "1.865.82"
$ws.Range("E3").Value = "  -2.73%  "

$ws.Range("D4").Formula = "'1.001"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Formula = "'317.36"
$ws.Range("E5").Value = "  -2.29%  "

$ws.Range("D6").Formula = "'1.000"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").Formula = "'0.4365"
$ws.Range("E7").Value = "  -5.16%  "

$ws.Range("D8").Formula = "'0.3720"
$ws.Range("E8").Value = "  -2.80%  "

$ws.Range("D9").Formula = "'0.07482"
$ws.Range("E9").Value = "  -2.91%  "

$ws.Range("D10").Formula = "'0.9351"
$ws.Range("E10").Value = "  -4.67%  "

$ws.Range("D11").Formula = "'21.33"
$ws.Range("E11").Value = "  -3.96%  "

$ws.Range("D12").Value = "1.874.61"
$ws.Range("E12").Value = "  -3.15%  "

$ws.Range("D13").Formula = "'6.738"
$ws.Range("E13").Value = "  -3.28%  "

$ws.Range("D14").Formula = "'5.441"
$ws.Range("E14").Value = "  -4.39%  "

$ws.Range("D15").Formula = "'0.06840"
$ws.Range("E15").Value = "  -2.04%  "

$ws.Range("E16").Value = "  +0.08%  "

$ws.Range("D17").Formula = "'81.63"
$ws.Range("E17").Value = "  -3.10%  "

$ws.Range("D18").Formula = "'0.000009055"
$ws.Range("E18").Value = "  -4.30%  "

$ws.Range("D19").Formula = "'0.9997"
$ws.Range("E19").Value = "  -0.09%  "

$ws.Range("D20").Formula = "'15.91"
$ws.Range("E20").Value = "  -4.64%  "

$ws.Range("D21").Value = "27.927.02"
$ws.Range("E21").Value = "  -3.61%  "

$ws.Range("D22").Formula = "'5.125"
$ws.Range("E22").Value = "  -3.99%  "

$ws.Range("D23").Formula = "'11.01"
$ws.Range("E23").Value = "  +0.57%  "

$ws.Range("D24").Value = "2.110.80"
$ws.Range("E24").Value = "  -2.15%  "

$ws.Range("D25").Formula = "'2.003"
$ws.Range("E25").Value = "  -4.25%  "

$ws.Range("D26").Formula = "'154.27"
$ws.Range("E26").Value = "  -2.61%  "

$ws.Range("D27").Formula = "'18.46"
$ws.Range("E27").Value = "  -3.07%  "

$ws.Range("D28").Formula = "'5.482"
$ws.Range("E28").Value = "  -3.89%  "

$ws.Range("D29").Formula = "'113.16"
$ws.Range("E29").Value = "  -3.95%  "

$ws.Range("D30").Formula = "'1.717"
$ws.Range("E30").Value = "  -8.03%  "

$ws.Range("D31").Formula = "'0.09013"
$ws.Range("E31").Value = "  -3.23%  "

$ws.Range("D32").Formula = "'0.8222"
$ws.Range("E32").Value = "  -5.13%  "

$ws.Range("D33").Formula = "'4.819"
$ws.Range("E33").Value = "  -5.78%  "

$ws.Range("E34").Value = "  -5.90%  "

$ws.Range("E35").Value = "  -2.24%  "

$ws.Range("D36").Formula = "'1.001"
$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("D37").Formula = "'0.05503"
$ws.Range("E37").Value = "  -3.56%  "

$ws.Range("D38").Formula = "'1.119"
$ws.Range("E38").Value = "  -3.27%  "

$ws.Range("D39").Formula = "'0.01979"
$ws.Range("E39").Value = "  -3.09%  "

$ws.Range("D40").Formula = "'2.965"
$ws.Range("E40").Value = "  -2.42%  "

$ws.Range("D41").Formula = "'0.5261"
$ws.Range("E41").Value = "  -4.50%  "

$ws.Range("D42").Formula = "'7.040"
$ws.Range("E42").Value = "  -6.49%  "

$ws.Range("D43").Formula = "'0.1702"
$ws.Range("E43").Value = "  -2.88%  "

$ws.Range("D44").Formula = "'8.785"
$ws.Range("E44").Value = "  -6.59%  "

$ws.Range("D45").Formula = "'0.06768"
$ws.Range("E45").Value = "  -2.04%  "

$ws.Range("D46").Formula = "'0.4899"
$ws.Range("E46").Value = "  -5.53%  "

$ws.Range("D47").Formula = "'10.67"
$ws.Range("E47").Value = "  -4.80%  "

$ws.Range("D48").Formula = "'107.30"
$ws.Range("E48").Value = "  -2.79%  "

$ws.Range("D49").Formula = "'1.679"
$ws.Range("E49").Value = "  -5.76%  "

$ws.Range("D50").Formula = "'0.9996"
$ws.Range("E50").Value = "  -0.09%  "

$ws.Range("D51").Formula = "'1.883"
$ws.Range("E51").Value = "  -13.82%  "
